# Insert a new price-report row for "Vega Modelo de Temuco - Betarraga" at
# row 593 (a newer weekly observation). All existing rows 593-688 shift down
# to 594-689 automatically, keeping their data untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(593).Insert()

$ws.Range("A593").Value = 10
$ws.Range("B593").Value = "Vega Modelo de Temuco"
$ws.Range("C593").Value = "La Araucanía"
$ws.Range("D593").Value = 45218
$ws.Range("E593").Value = 9
$ws.Range("F593").Value = 100114014
$ws.Range("G593").Value = "Betarraga"
$ws.Range("H593").Value = "Sin especificar"
$ws.Range("I593").Value = "Primera"
$ws.Range("J593").Value = 80
$ws.Range("K593").Value = 10000
$ws.Range("L593").Value = 11000
$ws.Range("M593").Value = 10312
$ws.Range("N593").Value = "$/saco 25 kilos"
$ws.Range("O593").Value = "Provincia de Cautín"
$ws.Range("P593").Value = 412
$ws.Range("Q593").Value = 25
$ws.Range("R593").Value = "Hortaliza"
